$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Num Objects" values (column H) for specific rows
$ws.Range("H3").Value = 9
$ws.Range("H9").Value = 7
$ws.Range("H10").Value = 10
$ws.Range("H11").Value = 10
$ws.Range("H12").Value = 9
$ws.Range("H13").Value = 9
$ws.Range("H14").Value = 11
$ws.Range("H15").Value = 16
$ws.Range("H16").Value = 8
$ws.Range("H23").Value = 17
$ws.Range("H27").Value = 5

# Add a data bar conditional format on H3:H44 (orange, matches the rest of the sheet's style)
$cf = $ws.Range("H3:H44").FormatConditions.AddDatabar()
$cf.BarColor.Color = 2668287

# Update selection to reflect the new cell the editor ended up on, and scroll so row 1 (not 7) is the top
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("H28").Select()
